$wb = $excel.ActiveWorkbook
Write-Host $wb.Worksheets.Count
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
